$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$bytes = [System.Text.Encoding]::UTF8.GetBytes($ws.Name)
Write-Host ($bytes -join ",")
